# Update the build timestamp embedded in the "Version" / citation text and
# in the per-row "S" column of the data sheet.

$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet  = $wb.Worksheets.Item("Boundaries and methane sources")

# A2: "Version: mines - January 30 (built on ... EST)"
$a2 = [string]$aboutSheet.Range("A2").Value()
$aboutSheet.Range("A2").Value = $a2.Replace($oldStamp, $newStamp)

# A6: "Recommended Citation: ... version 'mines - January 30 (built on ... EST)'. ..."
$a6 = [string]$aboutSheet.Range("A6").Value()
$aboutSheet.Range("A6").Value = $a6.Replace($oldStamp, $newStamp)

# S2:S36 on the data sheet: "mines - January 30 (built on ... EST)"
for ($row = 2; $row -le 36; $row++) {
    $cell = $dataSheet.Cells.Item($row, 19)
    $cellVal = [string]$cell.Value()
    $cell.Value = $cellVal.Replace($oldStamp, $newStamp)
}
